$wb = $excel.ActiveWorkbook

# --- 1. Select B47 on the existing "Tabelle1" sheet (it will no longer be the active tab) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B47").Select()

# --- 2. Add the new "Sheet1" worksheet after "Tabelle1" ---
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Sheet1"

# --- 3. Write the cells that introduce brand-new shared strings first, in the exact
#        order they were authored, so the shared-string table append order matches. ---
$ws2.Range("G3").Value = '42 - 300 (Avg.: 147)'
$ws2.Range("G2").Value = '0.3 - 12 (Avg.: 5.4)'
$ws2.Range("G6").Value = 'land required, fertilization, processing, construction of bioenergy plants'
$ws2.Range("C7").Value = 'none (doesn''t block land use)'
$ws2.Range("G10").Value = 'medium, possibly high competition for land'
$ws2.Range("B11").Value = 'can be used with agroforestry'
$ws2.Range("C11").Value = 'improved soil quality, reduced land erosion'
$ws2.Range("G12").Value = 'possible high competition for land, threat for food security'
$ws2.Range("G11").Value = 'electricity production, displacement of fossil fuels'
$ws2.Range("E11").Value = 'phytoplankton can increase oxygen content of oceans'
$ws2.Range("D11").Value = 'can improve soil fertility, reduce ocean acidity'

# --- 4. Fill in the remaining cells (all reuse already-known shared strings) ---
$ws2.Range("A1").Value = 'Technology'
$ws2.Range("B1").Value = 'Afforestation / Reforestation'
$ws2.Range("C1").Value = 'Soil Sequestration'
$ws2.Range("D1").Value = 'Enhanced Mineralization'
$ws2.Range("E1").Value = 'Ocean Fertilization'
$ws2.Range("F1").Value = 'DAC'
$ws2.Range("G1").Value = 'BECCS'
$ws2.Range("A2").Value = 'Potential'
$ws2.Range("B2").Value = '1.2 - 10 (Avg.: 5.8)'
$ws2.Range("C2").Value = '1.2 - 3.57 (Avg.: 2.4)'
$ws2.Range("D2").Value = '2.5 - 10 (Avg.: 4.9)'
$ws2.Range("E2").Value = '0.3 - 5 (Avg.: 2)'
$ws2.Range("F2").Value = '1.2 - 15 (Avg.: 7.1)'
$ws2.Range("A3").Value = 'Cost (USD/t CO2)'
$ws2.Range("B3").Value = '1 - 494'
$ws2.Range("C3").Value = '10 - 100'
$ws2.Range("D3").Value = '24 - 600'
$ws2.Range("E3").Value = '20 - 457'
$ws2.Range("F3").Value = '60 - 1000'
$ws2.Range("A4").Value = 'CAPEX'
$ws2.Range("B4").Value = 'low - medium'
$ws2.Range("C4").Value = 'medium'
$ws2.Range("D4").Value = 'medium - high'
$ws2.Range("E4").Value = 'low - medium'
$ws2.Range("F4").Value = 'high'
$ws2.Range("G4").Value = 'medium - high'
$ws2.Range("A5").Value = 'OPEX'
$ws2.Range("B5").Value = 'low'
$ws2.Range("C5").Value = 'low'
$ws2.Range("D5").Value = 'high'
$ws2.Range("E5").Value = 'medium'
$ws2.Range("F5").Value = 'high'
$ws2.Range("G5").Value = 'medium'
$ws2.Range("A6").Value = 'Cost drivers'
$ws2.Range("B6").Value = 'land required, management cost'
$ws2.Range("C6").Value = 'cost of adapting to new land management techniques'
$ws2.Range("D6").Value = 'construction of infrastructure, processing and transportation'
$ws2.Range("E6").Value = 'cost of mining and spreading nutrients'
$ws2.Range("F6").Value = 'construction of facilities, energy requirements'
$ws2.Range("A7").Value = 'Ressource requirements'
$ws2.Range("B7").Value = 'land, water'
$ws2.Range("D7").Value = 'rock, energy'
$ws2.Range("E7").Value = 'rock'
$ws2.Range("F7").Value = 'vast amounts of energy'
$ws2.Range("G7").Value = 'land, water, fertilizer'
$ws2.Range("A8").Value = 'Durability'
$ws2.Range("B8").Value = 'medium'
$ws2.Range("C8").Value = 'medium'
$ws2.Range("D8").Value = 'highest'
$ws2.Range("E8").Value = 'questionable'
$ws2.Range("F8").Value = 'depends on storage technology'
$ws2.Range("G8").Value = 'depends on storage technology'
$ws2.Range("A9").Value = 'Risks to durability'
$ws2.Range("B9").Value = 'fires, pests'
$ws2.Range("C9").Value = 'none, but requires continuous and permanent usage'
$ws2.Range("D9").Value = 'none'
$ws2.Range("E9").Value = 'none if sequesterd on ocean floor, but most co2 captured is respired back to surface quickly'
$ws2.Range("F9").Value = '-'
$ws2.Range("G9").Value = '-'
$ws2.Range("A10").Value = 'Additionality'
$ws2.Range("B10").Value = 'medium, converting farmland back into forests may result in forest removal in other locations'
$ws2.Range("C10").Value = 'high'
$ws2.Range("D10").Value = 'high'
$ws2.Range("E10").Value = 'questionable, due to possible nutrient robbing'
$ws2.Range("F10").Value = 'highest'
$ws2.Range("A11").Value = 'Co-Benefits'
$ws2.Range("F11").Value = 'can be used to clear air from pollution or draw water from the ambient air'
$ws2.Range("A12").Value = 'Negative Sideeffects'
$ws2.Range("B12").Value = 'possible competition for land'
$ws2.Range("C12").Value = 'none'
$ws2.Range("D12").Value = 'possible release of toxic metals to the food chain'
$ws2.Range("E12").Value = 'nutrient robbing, acidification of deep ocean'
$ws2.Range("F12").Value = 'co2 depletion of local ecosystems'
$ws2.Range("A13").Value = 'Verfication'
$ws2.Range("B13").Value = 'Somewhat difficult , but possible based on forest area'
$ws2.Range("C13").Value = 'Difficult'
$ws2.Range("D13").Value = 'Difficult'
$ws2.Range("E13").Value = 'Very difficult, requires measuring of carbon content of deep ocean'
$ws2.Range("F13").Value = 'Easy'
$ws2.Range("G13").Value = 'Easy'

# --- 5. Bold the row-header column (A), matching style index 1 used elsewhere ---
$ws2.Range("A1:A13").Font.Bold = $true

# --- 6. Column widths matching the source workbook (COM ColumnWidth quantizes to 1/6 char) ---
$ws2.Columns.Item(1).ColumnWidth = 28.8776041666667
$ws2.Columns.Item(2).ColumnWidth = 28.1666666666667
$ws2.Columns.Item(3).ColumnWidth = 27.0221354166667
$ws2.Columns.Item(4).ColumnWidth = 32.4518229166667
$ws2.Columns.Item(5).ColumnWidth = 41.1666666666667
$ws2.Columns.Item(6).ColumnWidth = 29.7369791666667
$ws2.Columns.Item(7).ColumnWidth = 28.8776041666667

# --- 7. Activate the new sheet and select D39 (final UI state) ---
$ws2.Activate()
$ws2.Range("D39").Select()
